# CRY_Warning_Template.docx update
#  - Collapse the spell/grammar-check-split merge-field runs (the
#    "{", "d.Something", "}" triplets Word had fragmented with
#    <w:proofErr/> markers) back into single, clean runs.
#  - Rename the d.DairyTestDataLoadDate merge field to d.ReportedOnDate.
#  - Append a literal "%" after the {d.DairyTestCryoPercent} field.
#  - Tidy a couple of other proofErr-fragmented sentences
#    (Milkoscann instrument, "12 month" period, Telephone/Fax lines).

$d = $word.ActiveDocument

function Set-ExactText($doc, $range, $newText) {
    # Forces Word to rewrite the range as a single plain run: first
    # collapse the whole range down to one throw-away character, then
    # expand that single character back out to the desired text. This
    # clears out any <w:proofErr/> markers and merges the touched runs.
    $range.Text = "@"
    $shrunk = $doc.Range($range.Start, $range.Start + 1)
    $shrunk.Text = $newText
}

function Replace-InParagraph($doc, $paraIndex, $searchText, $newText) {
    $p = $doc.Paragraphs.Item($paraIndex)
    $scope = $doc.Range($p.Range.Start, $p.Range.End)
    $found = $scope.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found in paragraph ${paraIndex}: $searchText"
    }
    Set-ExactText $doc $scope $newText
}

# {d.CurrentDate}
Replace-InParagraph $d 9 "{d.CurrentDate}" "{d.CurrentDate}"

# {d.IRMA_Num}
Replace-InParagraph $d 11 "{d.IRMA_Num}" "{d.IRMA_Num}"

# {d.LicenceHolderCompany} (leaves the trailing tab runs untouched)
Replace-InParagraph $d 13 "{d.LicenceHolderCompany}" "{d.LicenceHolderCompany}"

# {d.MailingAddress}
Replace-InParagraph $d 14 "{d.MailingAddress}" "{d.MailingAddress}"

# {d.MailingCity}, {d.MailingProv}, {d.PostCode} (single line, rebuild whole line)
Replace-InParagraph $d 15 "{d.MailingCity}, {d.MailingProv}, {d.PostCode}" "{d.MailingCity}, {d.MailingProv}, {d.PostCode}"

# " {d.DairyTestDataLoadDate}." -> " {d.ReportedOnDate}." (field rename)
Replace-InParagraph $d 19 " {d.DairyTestDataLoadDate}." " {d.ReportedOnDate}."

# "determined with a Milkoscann 6000 (FOSS) instrument ... {d.DairyTestCryoPercent}." -> adds a literal "%"
Replace-InParagraph $d 21 "determined with a Milkoscann 6000 (FOSS) instrument" "determined with a Milkoscann 6000 (FOSS) instrument"
Replace-InParagraph $d 21 "{d.DairyTestCryoPercent}." "{d.DairyTestCryoPercent}%."

# "Subsequent infractions ... 12 month period ..."
Replace-InParagraph $d 23 "Subsequent infractions over the tolerance level of 3.7% within a 12 month period will result in a cash penalty being applied." "Subsequent infractions over the tolerance level of 3.7% within a 12 month period will result in a cash penalty being applied."

# "Telephone:   (778) 666-0560"
Replace-InParagraph $d 32 "Telephone:   (778) 666-0560" "Telephone:   (778) 666-0560"

# "Fax:              (604) 556-3015"
Replace-InParagraph $d 33 "Fax:              (604) 556-3015" "Fax:              (604) 556-3015"
